$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (values updated slightly)
$ws.Range("B2").Value = 0.8325521762100273
$ws.Range("C2").Value = 0.8325521762100273
$ws.Range("D2").Value = 0.8325521762100273

# Row 3 - RandomForestRegressor (values updated)
$ws.Range("B3").Value = 0.9980834521943128
$ws.Range("C3").Value = 0.9980836383266043
$ws.Range("D3").Value = 0.9827637090982152

# Row 4 - renamed from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9975185867164996
$ws.Range("C4").Value = 0.997455495858747
$ws.Range("D4").Value = 0.9638017149380271

# Row 5 - renamed from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9989544420029824
$ws.Range("C5").Value = 0.9988894348326443
$ws.Range("D5").Value = 0.9981113888451142
